$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: image (col B), word (col C), category (col D)
# index (col A) is unchanged from the original workbook.
$stimuli = @(
    @("flower/flower101.png", "angeln", "flower"),
    @("dog/dog106.png", "lassen", "dog"),
    @("dog/dog102.png", "trotzen", "dog"),
    @("flower/flower099.png", "prüfen", "flower"),
    @("flower/flower097.png", "leeren", "flower"),
    @("dog/dog084.png", "bauen", "dog"),
    @("dog/dog103.png", "reisen", "dog"),
    @("dog/dog071.png", "deuten", "dog"),
    @("flower/flower098.png", "frischen", "flower"),
    @("dog/dog075.png", "heißen", "dog"),
    @("flower/flower095.png", "betteln", "flower"),
    @("flower/flower074.png", "piepen", "flower"),
    @("flower/flower070.png", "kennen", "flower"),
    @("flower/flower126.png", "tollen", "flower"),
    @("dog/dog104.png", "küssen", "dog"),
    @("flower/flower088.png", "legen", "flower"),
    @("flower/flower075.png", "lächeln", "flower"),
    @("flower/flower091.png", "mögen", "flower"),
    @("dog/dog115.png", "wählen", "dog"),
    @("dog/dog105.png", "holen", "dog"),
    @("flower/flower094.png", "wachsen", "flower"),
    @("flower/flower087.png", "sparen", "flower"),
    @("dog/dog107.png", "ändern", "dog"),
    @("flower/flower067.png", "rufen", "flower"),
    @("flower/flower080.png", "ärgern", "flower"),
    @("dog/dog085.png", "heben", "dog"),
    @("dog/dog073.png", "öffnen", "dog"),
    @("dog/dog072.png", "planen", "dog"),
    @("dog/dog114.png", "zögern", "dog"),
    @("dog/dog081.png", "narren", "dog"),
    @("dog/dog092.png", "hacken", "dog"),
    @("flower/flower104.png", "meinen", "flower")
)

$startRow = 2
for ($i = 0; $i -lt $stimuli.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $stimuli[$i][0]
    $ws.Cells.Item($r, 3).Value = $stimuli[$i][1]
    $ws.Cells.Item($r, 4).Value = $stimuli[$i][2]
}
